# Bulk Upload Template.xlsx - add 5 new employee rows (E0119-E0123 / p19-p23)
# to the Employees sheet, leaving Mobile/Location blank for the last couple of
# rows (optional fields) and moving the active selection/cursor down to the
# new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Employees")
$ws.Activate()

# Copy the formatting of the last existing data row (row 7) down into the
# five new rows (8-12) so the new cells pick up the same styles/number
# formats as the rest of the table.
$ws.Range("A7:K7").Copy()
$ws.Range("A8:K12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 8 - E0119 / p19
$ws.Range("A8").Value2 = 7
$ws.Range("B8").Value2 = "E0119"
$ws.Range("C8").Value2 = "p19"
$ws.Range("D8").Value2 = 44922
$ws.Range("E8").Value2 = 35967
$ws.Range("F8").Value2 = "Java Developer"
$ws.Range("G8").Value2 = "p19@gmail.com"
$ws.Range("H8").Value2 = "Female"
$ws.Range("I8").Value2 = 1234543268
$ws.Range("J8").Value2 = "Pune"
$ws.Range("K8").Value2 = 1

# Row 9 - E0120 / p20
$ws.Range("A9").Value2 = 8
$ws.Range("B9").Value2 = "E0120"
$ws.Range("C9").Value2 = "p20"
$ws.Range("D9").Value2 = 44922
$ws.Range("E9").Value2 = 35959
$ws.Range("F9").Value2 = "Solution Developer"
$ws.Range("G9").Value2 = "p20@gmail.com"
$ws.Range("H9").Value2 = "Female"
$ws.Range("I9").Value2 = 1234543268
$ws.Range("J9").Value2 = "Pune"
$ws.Range("K9").Value2 = 1

# Row 10 - E0121 / p21
$ws.Range("A10").Value2 = 9
$ws.Range("B10").Value2 = "E0121"
$ws.Range("C10").Value2 = "p21"
$ws.Range("D10").Value2 = 44922
$ws.Range("E10").Value2 = 35964
$ws.Range("F10").Value2 = "Java Developer"
$ws.Range("G10").Value2 = "p21@gmail.com"
$ws.Range("H10").Value2 = "Female"
$ws.Range("I10").Value2 = 1234543268
$ws.Range("J10").Value2 = "Pune"
$ws.Range("K10").Value2 = 1

# Row 11 - E0122 / p22 (Mobile left blank - optional field)
$ws.Range("A11").Value2 = 10
$ws.Range("B11").Value2 = "E0122"
$ws.Range("C11").Value2 = "p22"
$ws.Range("D11").Value2 = 44922
$ws.Range("E11").Value2 = 35964
$ws.Range("F11").Value2 = "Solution Developer"
$ws.Range("G11").Value2 = "p22@gmail.com"
$ws.Range("H11").Value2 = "Female"
$ws.Range("J11").Value2 = "Pune"
$ws.Range("K11").Value2 = 1

# Row 12 - E0123 / p23 (Mobile and Location left blank - optional fields)
$ws.Range("A12").Value2 = 11
$ws.Range("B12").Value2 = "E0123"
$ws.Range("C12").Value2 = "p23"
$ws.Range("D12").Value2 = 44922
$ws.Range("E12").Value2 = 35964
$ws.Range("F12").Value2 = "Solution Developer"
$ws.Range("G12").Value2 = "p23@gmail.com"
$ws.Range("H12").Value2 = "Female"
$ws.Range("K12").Value2 = 1

# Move the cursor / selection to the new last row, matching the place the
# author left off after the bulk upload.
$excel.Goto($ws.Range("G12"), $true)
